# Updated cryptos list data refresh
# Force the Price/Volume columns to Text format before writing so that
# Excel does not reinterpret numeric-looking strings (e.g. "578.59") as
# real numbers; the source data stores these as plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.116.08"
$ws.Range("E2").Value = "  -2.87%  "
$ws.Range("D3").Value = "3.518.32"
$ws.Range("E3").Value = "  -4.64%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "578.59"
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").Value = "171.57"
$ws.Range("E6").Value = "  -3.34%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.513.19"
$ws.Range("E7").Value = "  -4.46%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -5.34%  "
$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "0.582"
$ws.Range("E12").Value = "  -4.24%  "
$ws.Range("D13").Value = "46.99"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("D15").Value = "4.082.05"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("D16").Value = "8.54"
$ws.Range("E16").Value = "  -4.99%  "
$ws.Range("D17").Value = "622.57"
$ws.Range("E17").Value = "  -8.08%  "
$ws.Range("D18").Value = "3.516.90"
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").Value = "69.088.47"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "11.17"
$ws.Range("E22").Value = "  -3.46%  "
$ws.Range("D23").Value = "0.884"
$ws.Range("E23").Value = "  -5.90%  "
$ws.Range("D24").Value = "15.95"
$ws.Range("E24").Value = "  -7.99%  "
$ws.Range("D25").Value = "97.54"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -6.29%  "
$ws.Range("D29").Value = "9.35"
$ws.Range("E29").Value = "  -8.80%  "
$ws.Range("D30").Value = "32.64"
$ws.Range("E30").Value = "  -6.86%  "
$ws.Range("D31").Value = "3.17"
$ws.Range("E31").Value = "  -7.30%  "
$ws.Range("D32").Value = "8.55"
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("E33").Value = "  -7.09%  "
$ws.Range("D34").Value = "7.00"
$ws.Range("E34").Value = "  -6.38%  "
$ws.Range("D35").Value = "633.19"
$ws.Range("E35").Value = "  +8.01%  "
$ws.Range("D36").Value = "10.75"
$ws.Range("E36").Value = "  -3.80%  "
$ws.Range("E37").Value = "  -5.17%  "
$ws.Range("D38").Value = "3.42"
$ws.Range("E38").Value = "  -16.32%  "
$ws.Range("D39").Value = "56.72"
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("D41").Value = "0.0447"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  -5.65%  "
$ws.Range("D43").Value = "3.373.16"
$ws.Range("E43").Value = "  -8.35%  "
$ws.Range("E44").Value = "  -5.95%  "
$ws.Range("D45").Value = "32.93"
$ws.Range("E45").Value = "  -6.83%  "
$ws.Range("D46").Value = "0.0₃0691"
$ws.Range("E46").Value = "  -9.60%  "
$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -6.68%  "
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  +14.48%  "

# Restore the original (default/"Normal") cell style now that the text
# values are safely stored, so no stray number-format styling is left
# applied to the cells.
$dataRange.Style = "Normal"
